$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet is protected; unprotect to allow edits, re-protect afterwards
$ws.Unprotect()

# Update the confidential disclosure date from 2021-04-29 to 2021-04-30
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-04-30 for illustrative purposes only and are subject to change."

# Update Weight (column D) and Percent Change (column E) values for each holding row
$ws.Range("D2").Value = 0.01104855922273377
$ws.Range("E2").Value = -0.01156677181913779
$ws.Range("D3").Value = 0.01024692874285088
$ws.Range("E3").Value = -0.04237123420796896
$ws.Range("D4").Value = 0.01073354997825805
$ws.Range("E4").Value = -0.007916872835230016
$ws.Range("D5").Value = 0.01134652897874822
$ws.Range("E5").Value = -0.003929866989117237
$ws.Range("D6").Value = 0.01075744952051543
$ws.Range("E6").Value = -0.0007405578869414953
$ws.Range("D7").Value = 0.01204404154537096
$ws.Range("E7").Value = -0.001773049645390157
$ws.Range("D8").Value = 0.01118697740497442
$ws.Range("E8").Value = 0.008110300081102917
$ws.Range("D9").Value = 0.01125668440322511
$ws.Range("E9").Value = -0.01635605881890378
$ws.Range("D10").Value = 0.01040825065308819
$ws.Range("E10").Value = -0.003061615001913487
$ws.Range("D11").Value = 0.01100872665230481
$ws.Range("E11").Value = -0.01191014623850428
$ws.Range("D12").Value = 0.447897338188481
$ws.Range("E12").Value = -0.002549394518801873
$ws.Range("D13").Value = 0.01186147559757154
$ws.Range("E13").Value = 0.001585791309863582
$ws.Range("D14").Value = 0.01082317326172322
$ws.Range("E14").Value = 0.007013024187777273
$ws.Range("D15").Value = 0.01036974583500685
$ws.Range("E15").Value = 0.00610328638497637
$ws.Range("D16").Value = 0.01020112128685757
$ws.Range("E16").Value = -0.01156232374506494
$ws.Range("D17").Value = 0.01022369307676732
$ws.Range("E17").Value = -0.018452380952381
$ws.Range("D18").Value = 0.00854054568408567
$ws.Range("E18").Value = -0.05542312276519679
$ws.Range("D19").Value = 0.00880333000288786
$ws.Range("E19").Value = -0.02101479330844724
$ws.Range("D20").Value = 0.01239014232398706
$ws.Range("E20").Value = 0.01535988569387392
$ws.Range("D21").Value = 0.01193582972903894
$ws.Range("E21").Value = -0.02119138995494751
$ws.Range("D22").Value = 0.01185638588023895
$ws.Range("E22").Value = -0.03156146179401986
$ws.Range("D23").Value = 0.01185018970261667
$ws.Range("E23").Value = -0.0113445378151259
$ws.Range("D24").Value = 0.01183912509971973
$ws.Range("E24").Value = -0.00570093457943921
$ws.Range("D25").Value = 0.01234101548712467
$ws.Range("E25").Value = 0.000143451441687148
$ws.Range("D26").Value = 0.01130979449713039
$ws.Range("E26").Value = -0.02688424512796428
$ws.Range("D27").Value = 0.01034739533715504
$ws.Range("E27").Value = -0.01278898180029509
$ws.Range("D28").Value = 0.01246195159678817
$ws.Range("E28").Value = -0.03047172575446833
$ws.Range("D29").Value = 0.01023697060024364
$ws.Range("E29").Value = -0.006160830090791269
$ws.Range("D30").Value = 0.006992054508659709
$ws.Range("E30").Value = -0.02777206336144822
$ws.Range("D31").Value = 0.005309349700093937
$ws.Range("E31").Value = -0.01167031363967908
$ws.Range("D32").Value = 0.008990764375961928
$ws.Range("E32").Value = 0.002190580503833583
$ws.Range("D33").Value = 0.01070113069177004
$ws.Range("E33").Value = -0.002098950524737742
$ws.Range("D34").Value = 0.01041610652114501
$ws.Range("E34").Value = -0.002443195699975642
$ws.Range("D35").Value = 0.009555169769734547
$ws.Range("E35").Value = -0.01885175664095984
$ws.Range("D36").Value = 0.01099666623514715
$ws.Range("E36").Value = -0.07860262008733632
$ws.Range("D37").Value = 0.009914658717855945
$ws.Range("E37").Value = -0.009675583380762798
$ws.Range("D38").Value = 0.01146514152180335
$ws.Range("E38").Value = -0.01326963906581746
$ws.Range("D39").Value = 0.01332355222437244
$ws.Range("E39").Value = -0.007706617060855625
$ws.Range("D40").Value = 0.01122614609922957
$ws.Range("E40").Value = -0.01048689138576775
$ws.Range("D41").Value = 0.01173733075306794
$ws.Range("E41").Value = -0.008748114630467541
$ws.Range("D42").Value = 0.01122614609922957
$ws.Range("E42").Value = -0.0009363295880150391
$ws.Range("D43").Value = 0.01108462982817778
$ws.Range("E43").Value = 0
$ws.Range("D44").Value = 0.01106250062238391
$ws.Range("E44").Value = -0.02323441453876229
$ws.Range("D45").Value = 0.01177273748233813
$ws.Range("E45").Value = -0.006109022556391008
$ws.Range("D46").Value = 0.01071783824214441
$ws.Range("E46").Value = -0.0002271178741767299
$ws.Range("D47").Value = 0.01042351980508595
$ws.Range("E47").Value = -0.02674988854213112
$ws.Range("D48").Value = 0.01088756925058338
$ws.Range("E48").Value = -0.04634146341463397
$ws.Range("D49").Value = 0.009924948798550092
$ws.Range("E49").Value = -0.02307692307692311
$ws.Range("D50").Value = 0.00954786713182257
$ws.Range("E50").Value = 0.001251564455569509
$ws.Range("D51").Value = 0.009798812325525041
$ws.Range("E51").Value = -0.04861111111111105
$ws.Range("D52").Value = 0.01006878863621024
$ws.Range("E52").Value = -0.01252747252747255
$ws.Range("D53").Value = 0.009292053512845448
$ws.Range("E53").Value = -0.02914979757085012
$ws.Range("D54").Value = 0.004165822990695775
$ws.Range("E54").Value = -0.01394422310756971
$ws.Range("D55").Value = 0.004071773866071833
$ws.Range("E55").Value = 0.0004347826086956719
$ws.Range("D56").Value = 0.9999999999999997
$ws.Range("E56").Value = -0.008614235939379156

# Restore sheet protection
$ws.Protect()
